$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 1 de Octubre de 2020 a las 20:20"

# Row 4
$ws.Cells.Item(4, 2).Value = 7471839
$ws.Cells.Item(4, 3).Value = 24557
$ws.Cells.Item(4, 4).Value = 4716143
$ws.Cells.Item(4, 5).Value = 2543512
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 444
$ws.Cells.Item(4, 8).Value = 212184

# Row 14
$ws.Cells.Item(14, 2).Value = 577505
$ws.Cells.Item(14, 3).Value = 13970
$ws.Cells.Item(14, 4).Value = 96797
$ws.Cells.Item(14, 5).Value = 448689
$ws.Cells.Item(14, 6).Value = 0
$ws.Cells.Item(14, 7).Value = 63
$ws.Cells.Item(14, 8).Value = 32019

# Row 15
$ws.Cells.Item(15, 2).Value = 464750
$ws.Cells.Item(15, 3).Value = 1759
$ws.Cells.Item(15, 4).Value = 438148
$ws.Cells.Item(15, 5).Value = 13780
$ws.Cells.Item(15, 6).Value = 0
$ws.Cells.Item(15, 7).Value = 81
$ws.Cells.Item(15, 8).Value = 12822

# Row 25
$ws.Cells.Item(25, 2).Value = 294834
$ws.Cells.Item(25, 3).Value = 1923
$ws.Cells.Item(25, 4).Value = 257900
$ws.Cells.Item(25, 5).Value = 27353
$ws.Cells.Item(25, 6).Value = 0
$ws.Cells.Item(25, 7).Value = 10
$ws.Cells.Item(25, 8).Value = 9581

# Row 27
$ws.Cells.Item(27, 2).Value = 252533
$ws.Cells.Item(27, 3).Value = 7039
$ws.Cells.Item(27, 4).Value = 178976
$ws.Cells.Item(27, 5).Value = 71957
$ws.Cells.Item(27, 6).Value = 0
$ws.Cells.Item(27, 7).Value = 31
$ws.Cells.Item(27, 8).Value = 1600

# Row 33 (country -> Marruecos)
$ws.Cells.Item(33, 1).Value = "Marruecos"
$ws.Cells.Item(33, 2).Value = 126044
$ws.Cells.Item(33, 3).Value = 2391
$ws.Cells.Item(33, 4).Value = 104136
$ws.Cells.Item(33, 5).Value = 19679
$ws.Cells.Item(33, 6).Value = 0
$ws.Cells.Item(33, 7).Value = 35
$ws.Cells.Item(33, 8).Value = 2229

# Row 34 (country -> Catar)
$ws.Cells.Item(34, 1).Value = "Catar"
$ws.Cells.Item(34, 2).Value = 125959
$ws.Cells.Item(34, 3).Value = 199
$ws.Cells.Item(34, 4).Value = 122911
$ws.Cells.Item(34, 5).Value = 2834
$ws.Cells.Item(34, 6).Value = 0
$ws.Cells.Item(34, 7).Value = 0
$ws.Cells.Item(34, 8).Value = 214

# Row 35 (country -> Paises Bajos)
$ws.Cells.Item(35, 1).Value = "Paises Bajos"
$ws.Cells.Item(35, 2).Value = 124097
$ws.Cells.Item(35, 3).Value = 3252
$ws.Cells.Item(35, 4).Value = 0
$ws.Cells.Item(35, 5).Value = 0
$ws.Cells.Item(35, 6).Value = 0
$ws.Cells.Item(35, 7).Value = 13
$ws.Cells.Item(35, 8).Value = 6419

# Row 53 (country -> Etiopia)
$ws.Cells.Item(53, 1).Value = "Etiopia"
$ws.Cells.Item(53, 2).Value = 76098
$ws.Cells.Item(53, 3).Value = 730
$ws.Cells.Item(53, 4).Value = 31430
$ws.Cells.Item(53, 5).Value = 43463
$ws.Cells.Item(53, 6).Value = 0
$ws.Cells.Item(53, 7).Value = 7
$ws.Cells.Item(53, 8).Value = 1205

# Row 54 (country -> Costa Rica)
$ws.Cells.Item(54, 1).Value = "Costa Rica"
$ws.Cells.Item(54, 2).Value = 75760
$ws.Cells.Item(54, 3).Value = 0
$ws.Cells.Item(54, 4).Value = 37841
$ws.Cells.Item(54, 5).Value = 37015
$ws.Cells.Item(54, 6).Value = 0
$ws.Cells.Item(54, 7).Value = 0
$ws.Cells.Item(54, 8).Value = 904

# Row 74
$ws.Cells.Item(74, 2).Value = 36597
$ws.Cells.Item(74, 3).Value = 442
$ws.Cells.Item(74, 4).Value = 23364
$ws.Cells.Item(74, 5).Value = 11427
$ws.Cells.Item(74, 6).Value = 0
$ws.Cells.Item(74, 7).Value = 2
$ws.Cells.Item(74, 8).Value = 1806

# Row 106
$ws.Cells.Item(106, 2).Value = 9966
$ws.Cells.Item(106, 3).Value = 11
$ws.Cells.Item(106, 4).Value = 9613
$ws.Cells.Item(106, 5).Value = 286
$ws.Cells.Item(106, 6).Value = 0
$ws.Cells.Item(106, 7).Value = 1
$ws.Cells.Item(106, 8).Value = 67

# Row 108 (country -> Mozambique)
$ws.Cells.Item(108, 1).Value = "Mozambique"
$ws.Cells.Item(108, 2).Value = 8888
$ws.Cells.Item(108, 3).Value = 160
$ws.Cells.Item(108, 4).Value = 5573
$ws.Cells.Item(108, 5).Value = 3253
$ws.Cells.Item(108, 6).Value = 0
$ws.Cells.Item(108, 7).Value = 1
$ws.Cells.Item(108, 8).Value = 62

# Row 109 (country -> Haiti)
$ws.Cells.Item(109, 1).Value = "Haiti"
$ws.Cells.Item(109, 2).Value = 8766
$ws.Cells.Item(109, 3).Value = 0
$ws.Cells.Item(109, 4).Value = 6829
$ws.Cells.Item(109, 5).Value = 1708
$ws.Cells.Item(109, 6).Value = 0
$ws.Cells.Item(109, 7).Value = 0
$ws.Cells.Item(109, 8).Value = 229

# Row 110 (country -> Gabon)
$ws.Cells.Item(110, 1).Value = "Gabon"
$ws.Cells.Item(110, 2).Value = 8766
$ws.Cells.Item(110, 3).Value = 0
$ws.Cells.Item(110, 4).Value = 8005
$ws.Cells.Item(110, 5).Value = 707
$ws.Cells.Item(110, 6).Value = 0
$ws.Cells.Item(110, 7).Value = 0
$ws.Cells.Item(110, 8).Value = 54

# Row 119
$ws.Cells.Item(119, 2).Value = 5779
$ws.Cells.Item(119, 3).Value = 6
$ws.Cells.Item(119, 4).Value = 4514
$ws.Cells.Item(119, 5).Value = 1086
$ws.Cells.Item(119, 6).Value = 0
$ws.Cells.Item(119, 7).Value = 0
$ws.Cells.Item(119, 8).Value = 179

# Row 122
$ws.Cells.Item(122, 2).Value = 5500
$ws.Cells.Item(122, 3).Value = 18
$ws.Cells.Item(122, 4).Value = 5000
$ws.Cells.Item(122, 5).Value = 389
$ws.Cells.Item(122, 6).Value = 0
$ws.Cells.Item(122, 7).Value = 2
$ws.Cells.Item(122, 8).Value = 111

# Row 136 (country -> Aruba)
$ws.Cells.Item(136, 1).Value = "Aruba"
$ws.Cells.Item(136, 2).Value = 3998
$ws.Cells.Item(136, 3).Value = 35
$ws.Cells.Item(136, 4).Value = 3327
$ws.Cells.Item(136, 5).Value = 644
$ws.Cells.Item(136, 6).Value = 0
$ws.Cells.Item(136, 7).Value = 0
$ws.Cells.Item(136, 8).Value = 27

# Row 137 (country -> Reunion)
$ws.Cells.Item(137, 1).Value = "Reunion"
$ws.Cells.Item(137, 2).Value = 3993
$ws.Cells.Item(137, 3).Value = 0
$ws.Cells.Item(137, 4).Value = 2819
$ws.Cells.Item(137, 5).Value = 1158
$ws.Cells.Item(137, 6).Value = 0
$ws.Cells.Item(137, 7).Value = 0
$ws.Cells.Item(137, 8).Value = 16

# Row 152
$ws.Cells.Item(152, 2).Value = 2238
$ws.Cells.Item(152, 3).Value = 7
$ws.Cells.Item(152, 4).Value = 1695
$ws.Cells.Item(152, 5).Value = 471
$ws.Cells.Item(152, 6).Value = 0
$ws.Cells.Item(152, 7).Value = 0
$ws.Cells.Item(152, 8).Value = 72

# Row 176 (country -> Burundi)
$ws.Cells.Item(176, 1).Value = "Burundi"
$ws.Cells.Item(176, 2).Value = 510
$ws.Cells.Item(176, 3).Value = 2
$ws.Cells.Item(176, 4).Value = 472
$ws.Cells.Item(176, 5).Value = 37
$ws.Cells.Item(176, 6).Value = 0
$ws.Cells.Item(176, 7).Value = 0
$ws.Cells.Item(176, 8).Value = 1

# Row 177 (country -> Tanzania)
$ws.Cells.Item(177, 1).Value = "Tanzania"
$ws.Cells.Item(177, 2).Value = 509
$ws.Cells.Item(177, 3).Value = 0
$ws.Cells.Item(177, 4).Value = 183
$ws.Cells.Item(177, 5).Value = 305
$ws.Cells.Item(177, 6).Value = 0
$ws.Cells.Item(177, 7).Value = 0
$ws.Cells.Item(177, 8).Value = 21

# Row 189
$ws.Cells.Item(189, 2).Value = 219
$ws.Cells.Item(189, 3).Value = 1
$ws.Cells.Item(189, 4).Value = 187
$ws.Cells.Item(189, 5).Value = 30
$ws.Cells.Item(189, 6).Value = 0
$ws.Cells.Item(189, 7).Value = 0
$ws.Cells.Item(189, 8).Value = 2

# Row 194
$ws.Cells.Item(194, 2).Value = 144
$ws.Cells.Item(194, 3).Value = 0
$ws.Cells.Item(194, 4).Value = 143
$ws.Cells.Item(194, 5).Value = 1
$ws.Cells.Item(194, 6).Value = 0
$ws.Cells.Item(194, 7).Value = 0
$ws.Cells.Item(194, 8).Value = 0

# Row 196
$ws.Cells.Item(196, 2).Value = 119
$ws.Cells.Item(196, 3).Value = 1
$ws.Cells.Item(196, 4).Value = 114
$ws.Cells.Item(196, 5).Value = 4
$ws.Cells.Item(196, 6).Value = 0
$ws.Cells.Item(196, 7).Value = 0
$ws.Cells.Item(196, 8).Value = 1

# Row 207 (country -> Nueva Caledonia)
$ws.Cells.Item(207, 1).Value = "Nueva Caledonia"
$ws.Cells.Item(207, 2).Value = 27
$ws.Cells.Item(207, 3).Value = 0
$ws.Cells.Item(207, 4).Value = 27
$ws.Cells.Item(207, 5).Value = 0
$ws.Cells.Item(207, 6).Value = 0
$ws.Cells.Item(207, 7).Value = 0
$ws.Cells.Item(207, 8).Value = 0

# Row 208 (country -> Santa Lucia)
$ws.Cells.Item(208, 1).Value = "Santa Lucia"
$ws.Cells.Item(208, 2).Value = 27
$ws.Cells.Item(208, 3).Value = 0
$ws.Cells.Item(208, 4).Value = 27
$ws.Cells.Item(208, 5).Value = 0
$ws.Cells.Item(208, 6).Value = 0
$ws.Cells.Item(208, 7).Value = 0
$ws.Cells.Item(208, 8).Value = 0
